$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. D-column price values are forced to text
# (NumberFormat "@") before assignment so numeric-looking strings like
# "356.33" are stored as text instead of being coerced into numbers,
# matching the original inlineStr cell type. Style is reset to "Normal"
# afterwards so no stray cell-style index is left behind.
$updates = @(
    @{ Cell = "D2"; Value = '51.790.27'; Text = $true }
    @{ Cell = "E2"; Value = '  -0.26%  '; Text = $false }
    @{ Cell = "D3"; Value = '2.778.88'; Text = $true }
    @{ Cell = "E3"; Value = '  -1.21%  '; Text = $false }
    @{ Cell = "E4"; Value = '  +0.03%  '; Text = $false }
    @{ Cell = "D5"; Value = '356.33'; Text = $true }
    @{ Cell = "E5"; Value = '  +0.98%  '; Text = $false }
    @{ Cell = "D6"; Value = '108.93'; Text = $true }
    @{ Cell = "E6"; Value = '  -2.20%  '; Text = $false }
    @{ Cell = "E7"; Value = '  -1.53%  '; Text = $false }
    @{ Cell = "D8"; Value = '0.999'; Text = $true }
    @{ Cell = "E8"; Value = '  +0.05%  '; Text = $false }
    @{ Cell = "D9"; Value = '0.586'; Text = $true }
    @{ Cell = "E9"; Value = '  -1.82%  '; Text = $false }
    @{ Cell = "D10"; Value = '39.68'; Text = $true }
    @{ Cell = "E10"; Value = '  -2.11%  '; Text = $false }
    @{ Cell = "E11"; Value = '  +2.64%  '; Text = $false }
    @{ Cell = "D12"; Value = '0.0844'; Text = $true }
    @{ Cell = "E12"; Value = '  -1.17%  '; Text = $false }
    @{ Cell = "D13"; Value = '19.48'; Text = $true }
    @{ Cell = "E13"; Value = '  -1.56%  '; Text = $false }
    @{ Cell = "D14"; Value = '7.61'; Text = $true }
    @{ Cell = "E14"; Value = '  -1.93%  '; Text = $false }
    @{ Cell = "D15"; Value = '3.211.09'; Text = $true }
    @{ Cell = "E15"; Value = '  -1.29%  '; Text = $false }
    @{ Cell = "D16"; Value = '2.790.30'; Text = $true }
    @{ Cell = "E16"; Value = '  -0.58%  '; Text = $false }
    @{ Cell = "D17"; Value = '0.933'; Text = $true }
    @{ Cell = "E17"; Value = '  +1.25%  '; Text = $false }
    @{ Cell = "D18"; Value = '51.681.42'; Text = $true }
    @{ Cell = "E18"; Value = '  +0.01%  '; Text = $false }
    @{ Cell = "E19"; Value = '  -0.95%  '; Text = $false }
    @{ Cell = "D20"; Value = '3.10'; Text = $true }
    @{ Cell = "E20"; Value = '  -0.18%  '; Text = $false }
    @{ Cell = "D21"; Value = '13.14'; Text = $true }
    @{ Cell = "E21"; Value = '  -1.21%  '; Text = $false }
    @{ Cell = "E22"; Value = '  -2.17%  '; Text = $false }
    @{ Cell = "D23"; Value = '70.17'; Text = $true }
    @{ Cell = "E23"; Value = '  -0.07%  '; Text = $false }
    @{ Cell = "D24"; Value = '269.09'; Text = $true }
    @{ Cell = "E24"; Value = '  +0.41%  '; Text = $false }
    @{ Cell = "E25"; Value = '  -2.43%  '; Text = $false }
    @{ Cell = "D26"; Value = '26.40'; Text = $true }
    @{ Cell = "E26"; Value = '  -1.64%  '; Text = $false }
    @{ Cell = "E27"; Value = '  +0.01%  '; Text = $false }
    @{ Cell = "E28"; Value = '  +16.48%  '; Text = $false }
    @{ Cell = "D29"; Value = '10.23'; Text = $true }
    @{ Cell = "E29"; Value = '  -0.19%  '; Text = $false }
    @{ Cell = "E30"; Value = '  -1.46%  '; Text = $false }
    @{ Cell = "D31"; Value = '6.26'; Text = $true }
    @{ Cell = "E31"; Value = '  +6.11%  '; Text = $false }
    @{ Cell = "D32"; Value = '34.88'; Text = $true }
    @{ Cell = "E32"; Value = '  +1.66%  '; Text = $false }
    @{ Cell = "D33"; Value = '51.63'; Text = $true }
    @{ Cell = "E33"; Value = '  -1.83%  '; Text = $false }
    @{ Cell = "D34"; Value = '0.0450'; Text = $true }
    @{ Cell = "E34"; Value = '  -8.91%  '; Text = $false }
    @{ Cell = "D35"; Value = '0.0838'; Text = $true }
    @{ Cell = "E35"; Value = '  -0.76%  '; Text = $false }
    @{ Cell = "D36"; Value = '5.11'; Text = $true }
    @{ Cell = "E36"; Value = '  -6.68%  '; Text = $false }
    @{ Cell = "D37"; Value = '0.999'; Text = $true }
    @{ Cell = "E37"; Value = '  +0.04%  '; Text = $false }
    @{ Cell = "D38"; Value = '18.65'; Text = $true }
    @{ Cell = "E38"; Value = '  +2.31%  '; Text = $false }
    @{ Cell = "D39"; Value = '3.14'; Text = $true }
    @{ Cell = "E39"; Value = '  -2.88%  '; Text = $false }
    @{ Cell = "D40"; Value = '1.96'; Text = $true }
    @{ Cell = "E40"; Value = '  -3.21%  '; Text = $false }
    @{ Cell = "E41"; Value = '  +2.92%  '; Text = $false }
    @{ Cell = "E42"; Value = '  -2.10%  '; Text = $false }
    @{ Cell = "E43"; Value = '  -2.32%  '; Text = $false }
    @{ Cell = "D44"; Value = '119.17'; Text = $true }
    @{ Cell = "E44"; Value = '  -6.22%  '; Text = $false }
    @{ Cell = "D45"; Value = '21.62'; Text = $true }
    @{ Cell = "E45"; Value = '  -5.78%  '; Text = $false }
    @{ Cell = "D46"; Value = '2.081.76'; Text = $true }
    @{ Cell = "E46"; Value = '  -0.21%  '; Text = $false }
    @{ Cell = "D47"; Value = '3.28'; Text = $true }
    @{ Cell = "E47"; Value = '  -1.36%  '; Text = $false }
    @{ Cell = "E48"; Value = '  +1.06%  '; Text = $false }
    @{ Cell = "D49"; Value = '0.944'; Text = $true }
    @{ Cell = "E49"; Value = '  -3.11%  '; Text = $false }
    @{ Cell = "E50"; Value = '  -5.98%  '; Text = $false }
    @{ Cell = "B51"; Value = 'MultiversX'; Text = $false }
    @{ Cell = "C51"; Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'; Text = $false }
    @{ Cell = "D51"; Value = '58.44'; Text = $true }
    @{ Cell = "E51"; Value = '  -3.55%  '; Text = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Text) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
